$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cryptocurrency price/volume data, including the two
# row swaps (Kaspa <-> Fetch.AI, Stacks <-> FirstDigitalUSD).
# Every assigned value starts with a literal leading apostrophe, which is
# Excel's "store as text" marker -- this keeps numeric-looking strings
# (e.g. "570.36") stored as literal text, matching the original inlineStr
# cells, instead of letting Excel auto-convert them into real numbers.
$ws.Range("D2").Value = "'60.824.54"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("D3").Value = "'3.375.48"
$ws.Range("E3").Value = "'  -0.21%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'570.36"
$ws.Range("E5").Value = "'  -0.63%  "
$ws.Range("D6").Value = "'136.26"
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("D8").Value = "'3.371.94"
$ws.Range("E8").Value = "'  -0.29%  "
$ws.Range("E9").Value = "'  -1.01%  "
$ws.Range("E10").Value = "'  +1.53%  "
$ws.Range("D11").Value = "'0.122"
$ws.Range("E11").Value = "'  -2.64%  "
$ws.Range("D12").Value = "'0.379"
$ws.Range("E12").Value = "'  -2.67%  "
$ws.Range("D13").Value = "'3.950.20"
$ws.Range("E13").Value = "'  -0.25%  "
$ws.Range("E14").Value = "'  -0.80%  "
$ws.Range("D15").Value = "'26.05"
$ws.Range("E15").Value = "'  +0.93%  "
$ws.Range("D16").Value = "'3.373.10"
$ws.Range("E16").Value = "'  -0.33%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "'  -3.45%  "
$ws.Range("D18").Value = "'60.902.39"
$ws.Range("E18").Value = "'  -1.00%  "
$ws.Range("D19").Value = "'13.79"
$ws.Range("E19").Value = "'  -2.03%  "
$ws.Range("E20").Value = "'  -0.99%  "
$ws.Range("D21").Value = "'9.21"
$ws.Range("E21").Value = "'  -1.61%  "
$ws.Range("D22").Value = "'371.67"
$ws.Range("E22").Value = "'  -1.03%  "
$ws.Range("D23").Value = "'3.506.96"
$ws.Range("E23").Value = "'  -0.51%  "
$ws.Range("D24").Value = "'0.548"
$ws.Range("E24").Value = "'  -1.56%  "
$ws.Range("E25").Value = "'  +0.12%  "
$ws.Range("D26").Value = "'70.63"
$ws.Range("E26").Value = "'  -0.82%  "
$ws.Range("D27").Value = "'0.0000122"
$ws.Range("E27").Value = "'  -2.61%  "
$ws.Range("B28").Value = "'Fetch.AI"
$ws.Range("C28").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.60"
$ws.Range("E28").Value = "'  -7.21%  "
$ws.Range("B29").Value = "'Kaspa"
$ws.Range("C29").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.172"
$ws.Range("E29").Value = "'  +6.67%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  -0.03%  "
$ws.Range("E31").Value = "'  -2.28%  "
$ws.Range("D32").Value = "'8.02"
$ws.Range("E32").Value = "'  -2.70%  "
$ws.Range("E33").Value = "'  -2.21%  "
$ws.Range("E34").Value = "'  -0.06%  "
$ws.Range("D35").Value = "'23.30"
$ws.Range("E35").Value = "'  -0.68%  "
$ws.Range("E36").Value = "'  -3.31%  "
$ws.Range("E37").Value = "'  -0.39%  "
$ws.Range("D38").Value = "'6.77"
$ws.Range("E38").Value = "'  -0.65%  "
$ws.Range("D39").Value = "'164.70"
$ws.Range("E39").Value = "'  -0.51%  "
$ws.Range("D40").Value = "'0.0762"
$ws.Range("E40").Value = "'  -1.83%  "
$ws.Range("D41").Value = "'25.54"
$ws.Range("E41").Value = "'  +3.20%  "
$ws.Range("B42").Value = "'Stacks"
$ws.Range("C42").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.74"
$ws.Range("E42").Value = "'  +1.06%  "
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("D44").Value = "'0.771"
$ws.Range("E44").Value = "'  -0.74%  "
$ws.Range("D45").Value = "'41.87"
$ws.Range("E45").Value = "'  +1.20%  "
$ws.Range("D46").Value = "'4.34"
$ws.Range("E46").Value = "'  -1.58%  "
$ws.Range("E47").Value = "'  -5.32%  "
$ws.Range("D48").Value = "'2.503.57"
$ws.Range("E48").Value = "'  +6.86%  "
$ws.Range("D49").Value = "'23.53"
$ws.Range("E49").Value = "'  +3.39%  "
$ws.Range("D50").Value = "'6.74"
$ws.Range("E50").Value = "'  -1.25%  "
$ws.Range("D51").Value = "'2.39"
$ws.Range("E51").Value = "'  +0.93%  "
